# PowerEnJoyGantt: translate diagram columns / sheet to English.
#
# Summary of the change (per the commit "translate diagrams columns to
# english"):
#   - sheet renamed RASD -> Sheet1
#   - the 9 column headers in row 1 are translated from Italian to English
#     (the task/resource/comment cell text underneath is already English
#     and is unchanged)
#   - the stray/leftover formatted row 25 (a single empty, styled D25 cell)
#     is removed, shrinking the used range from A1:V25 to A1:V19
#   - the selection marker follows the old D25 cell, which is now D24
#   - a couple of column widths are nudged to reflect the shorter English
#     header text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet tab rename -----------------------------------------------------
$ws.Name = "Sheet1"

# --- translate the header row (row 1) -------------------------------------
$ws.Range("A1").Value = "Task #"
$ws.Range("B1").Value = "Task Name"
$ws.Range("C1").Value = "Resource"
$ws.Range("D1").Value = "Task Notes"
$ws.Range("E1").Value = "Duration"
$ws.Range("F1").Value = "Start Date"
$ws.Range("G1").Value = "Finish Date"
$ws.Range("H1").Value = "Dependency"
$ws.Range("I1").Value = "Outline Level"

# --- drop the trailing empty styled row (old row 25) ----------------------
$ws.Rows.Item(25).Delete()

# Selection tracks the old D25 cell, now shifted up to D24 after the delete.
$ws.Range("D24").Select()

# --- column width touch-ups (narrower English labels) ---------------------
# Column "ColumnWidth" is expressed in characters; the engine stores the
# column width as (round(ColumnWidth*6)+5)/6, so these inputs reproduce the
# exact target widths of 7 / 13 / 13 / 14 characters for columns A, I, J, L.
$ws.Columns.Item(1).ColumnWidth = 6.166666666666667   # -> stored width 7
$ws.Columns.Item(9).ColumnWidth = 12.166666666666666  # -> stored width 13
$ws.Columns.Item(10).ColumnWidth = 12.166666666666666 # -> stored width 13
$ws.Columns.Item(12).ColumnWidth = 13.166666666666666 # -> stored width 14
